# Readme.docx edit: expand the "navigation bar" paragraph with inline
# <ul>/<div>/<id>/<li>/<.active>/<href> tag call-outs, add the "unordered
# list" phrase (carrying the _GoBack bookmark to its new location).
#
# Word's own find/replace can only work with literal text (the engine
# normalizes adjacent same-formatted runs back together on save, exactly
# like the underlying <w:r> splits in the original file), so the edit is
# expressed as one literal-text Find/Replace over the whole affected
# sentence span, followed by re-seating the singleton "_GoBack" bookmark
# (Word only ever keeps one; adding a new one retires the old one
# automatically, matching the diff's bookmark relocation).

$d = $word.ActiveDocument

$oldText = "used for the web page" +
    ". This took quite some time because it needed to be in several" +
    " div and id tags in the html." +
    " " +
    "Each page should be linked to the list" +
    " that is assigned" +
    " f" +
    "or the home" +
    " or the index page. The list that the home page falls under should not be navigated as that" + [char]0x2019 + "s the page you are in therefore to avoid this an active element was inserted" +
    " " +
    "though it is in a " +
    "href" +
    "."

$newText = "used for the web page" +
    " in an <ul> " +
    "unordered list" +
    ". This took quite some time because it needed to be in several" +
    " " +
    "<div> and <id> tags in the html." +
    " " +
    "Each page should be linked to the list " +
    "<li>" +
    "that is assigned" +
    " f" +
    "or the home or the index page. The list that the home page falls under should not be navigated as that" + [char]0x2019 + "s the page you are in therefore to avoid this an " +
    "<.active>" +
    " element was inserted though it is in a " +
    "<href>" +
    "."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# Re-seat the _GoBack bookmark right after "<ul> " / right before "unordered list",
# which is where the author's cursor ended up. Word keeps only one _GoBack
# bookmark, so this also removes it from its old spot (end of the document,
# just before "Page 4 - Contact").
$anchor = $d.Content
$anchor.Find.Execute("in an <ul> ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
$d.Bookmarks.Add("_GoBack", $anchor)
